$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 1282000
$ws.Range("E8").Value = 1216100
$ws.Range("F8").Value = 940300
$ws.Range("G8").Value = 726000
$ws.Range("H8").Value = 672500
$ws.Range("I8").Value = 604400
$ws.Range("J8").Value = 459200
$ws.Range("D9").Value = 999200
$ws.Range("E9").Value = 852600
$ws.Range("F9").Value = 629600
$ws.Range("G9").Value = 567400
$ws.Range("H9").Value = 520800
$ws.Range("I9").Value = 463700
$ws.Range("J9").Value = 367000
$ws.Range("D10").Value = 282800
$ws.Range("E10").Value = 363500
$ws.Range("F10").Value = 310700
$ws.Range("G10").Value = 158600
$ws.Range("H10").Value = 151600
$ws.Range("I10").Value = 140700
$ws.Range("J10").Value = 92200
$ws.Range("D17").Value = 1284000
$ws.Range("E17").Value = 1074400
$ws.Range("F17").Value = 810400
$ws.Range("G17").Value = 715400
$ws.Range("H17").Value = 656100
$ws.Range("I17").Value = 584900
$ws.Range("J17").Value = 478500
$ws.Range("D18").Value = -2000
$ws.Range("E18").Value = 141700
$ws.Range("F18").Value = 129800
$ws.Range("G18").Value = 10600
$ws.Range("H18").Value = 16400
$ws.Range("I18").Value = 19600
$ws.Range("J18").Value = -19300
$ws.Range("D20").Value = -38100
$ws.Range("E20").Value = 115800
$ws.Range("F20").Value = 51300
$ws.Range("G20").Value = 22700
$ws.Range("H20").Value = 2200
$ws.Range("I20").Value = -4200
$ws.Range("J20").Value = 6000
$ws.Range("D21").Value = -11800
$ws.Range("E21").Value = 285200
$ws.Range("F21").Value = 204700
$ws.Range("G21").Value = 51000
$ws.Range("H21").Value = 34200
$ws.Range("I21").Value = 26200
$ws.Range("J21").Value = "NA"
$ws.Range("D22").Value = 1900
$ws.Range("H22").Value = 4000
$ws.Range("I22").Value = 4600
$ws.Range("J22").Value = 3000
$ws.Range("D23").Value = -42000
$ws.Range("E23").Value = 257400
$ws.Range("F23").Value = 181100
$ws.Range("G23").Value = 33300
$ws.Range("H23").Value = 14600
$ws.Range("I23").Value = 10700
$ws.Range("J23").Value = -16300
$ws.Range("D24").Value = -8300
$ws.Range("E24").Value = 75400
$ws.Range("F24").Value = 53700
$ws.Range("G24").Value = 2000
$ws.Range("D26").Value = -33700
$ws.Range("E26").Value = 182000
$ws.Range("F26").Value = 127400
$ws.Range("G26").Value = 31300
$ws.Range("H26").Value = 13700
$ws.Range("I26").Value = 10500
$ws.Range("J26").Value = -16300
$ws.Range("D27").Value = -33700
$ws.Range("E27").Value = 182000
$ws.Range("F27").Value = 127400
$ws.Range("G27").Value = 31300
$ws.Range("H27").Value = 13900
$ws.Range("I27").Value = 11100
$ws.Range("J27").Value = -15200
$ws.Range("D32").Value = 38100
$ws.Range("E32").Value = -115800
$ws.Range("F32").Value = -51300
$ws.Range("G32").Value = -22700
$ws.Range("H32").Value = -2200
$ws.Range("I32").Value = 4200
$ws.Range("J32").Value = -6000
$ws.Range("D33").Value = -33700
$ws.Range("E33").Value = 182000
$ws.Range("F33").Value = 127400
$ws.Range("G33").Value = 31300
$ws.Range("H33").Value = 13900
$ws.Range("I33").Value = 11100
$ws.Range("J33").Value = -15200
$ws.Range("D35").Value = -33700
$ws.Range("E35").Value = 182000
$ws.Range("F35").Value = 127400
$ws.Range("G35").Value = 31300
$ws.Range("H35").Value = 13900
$ws.Range("I35").Value = 11100
$ws.Range("J35").Value = -15200
$ws.Range("D41").Value = 359500
$ws.Range("E41").Value = 365700
$ws.Range("F41").Value = 266700
$ws.Range("G41").Value = 234300
$ws.Range("H41").Value = 126800
$ws.Range("I41").Value = 42500
$ws.Range("J41").Value = 22800
$ws.Range("D42").Value = 25700
$ws.Range("E42").Value = 28100
$ws.Range("G42").Value = 3200
$ws.Range("D43").Value = 54200
$ws.Range("E43").Value = 32100
$ws.Range("F43").Value = 24800
$ws.Range("G43").Value = 35300
$ws.Range("H43").Value = 31100
$ws.Range("I43").Value = 39900
$ws.Range("J43").Value = 12400
$ws.Range("D44").Value = 15200
$ws.Range("E44").Value = 12600
$ws.Range("F44").Value = 8400
$ws.Range("G44").Value = 14400
$ws.Range("H44").Value = 5900
$ws.Range("I44").Value = 10000
$ws.Range("J44").Value = 4100
$ws.Range("D45").Value = 130400
$ws.Range("E45").Value = 158900
$ws.Range("F45").Value = 113800
$ws.Range("G45").Value = 40000
$ws.Range("H45").Value = 42500
$ws.Range("I45").Value = 40000
$ws.Range("J45").Value = 19100
$ws.Range("D46").Value = 585100
$ws.Range("E46").Value = 597400
$ws.Range("F46").Value = 373600
$ws.Range("G46").Value = 190800
$ws.Range("H46").Value = 206900
$ws.Range("I46").Value = 93900
$ws.Range("J46").Value = 58500
$ws.Range("D47").Value = 315400
$ws.Range("E47").Value = 356000
$ws.Range("F47").Value = 246200
$ws.Range("G47").Value = 183400
$ws.Range("H47").Value = 134700
$ws.Range("I47").Value = 116100
$ws.Range("J47").Value = 103500
$ws.Range("D48").Value = 226300
$ws.Range("E48").Value = 130600
$ws.Range("F48").Value = 131900
$ws.Range("G48").Value = 230000
$ws.Range("H48").Value = 69400
$ws.Range("I48").Value = 123600
$ws.Range("J48").Value = 78500
$ws.Range("D49").Value = 9800
$ws.Range("E49").Value = 5900
$ws.Range("F49").Value = 4900
$ws.Range("G49").Value = 7500
$ws.Range("H49").Value = 4100
$ws.Range("I49").Value = 6200
$ws.Range("J49").Value = 3000
$ws.Range("D52").Value = 35600
$ws.Range("E52").Value = 36600
$ws.Range("F52").Value = 31100
$ws.Range("G52").Value = 204900
$ws.Range("H52").Value = 18300
$ws.Range("I52").Value = 39300
$ws.Range("J52").Value = 18300
$ws.Range("D54").Value = 1172300
$ws.Range("E54").Value = 1126600
$ws.Range("F54").Value = 787800
$ws.Range("G54").Value = 512300
$ws.Range("H54").Value = 433300
$ws.Range("I54").Value = 294900
$ws.Range("J54").Value = 261800
$ws.Range("D57").Value = 57800
$ws.Range("E57").Value = 47900
$ws.Range("F57").Value = 41100
$ws.Range("G57").Value = 26200
$ws.Range("H57").Value = 27800
$ws.Range("I57").Value = 54200
$ws.Range("J57").Value = 27800
$ws.Range("D58").Value = 124300
$ws.Range("E58").Value = 54400
$ws.Range("F58").Value = 70900
$ws.Range("G58").Value = 42300
$ws.Range("H58").Value = 13900
$ws.Range("I58").Value = 27300
$ws.Range("J58").Value = 35500
$ws.Range("D59").Value = 302600
$ws.Range("E59").Value = 305700
$ws.Range("F59").Value = 253800
$ws.Range("G59").Value = 271200
$ws.Range("H59").Value = 158600
$ws.Range("I59").Value = 206700
$ws.Range("J59").Value = 101100
$ws.Range("D60").Value = 484700
$ws.Range("E60").Value = 408000
$ws.Range("F60").Value = 365800
$ws.Range("G60").Value = 246600
$ws.Range("H60").Value = 200200
$ws.Range("I60").Value = 192500
$ws.Range("J60").Value = 164400
$ws.Range("D61").Value = 55800
$ws.Range("E61").Value = 48800
$ws.Range("F61").Value = 11400
$ws.Range("G61").Value = 22000
$ws.Range("H61").Value = 19000
$ws.Range("I61").Value = 38400
$ws.Range("J61").Value = 44500
$ws.Range("D62").Value = 106100
$ws.Range("E62").Value = 111500
$ws.Range("F62").Value = 57600
$ws.Range("G62").Value = 13600
$ws.Range("H62").Value = 9100
$ws.Range("I62").Value = 8400
$ws.Range("J62").Value = 8500
$ws.Range("D66").Value = 646700
$ws.Range("E66").Value = 568300
$ws.Range("F66").Value = 434800
$ws.Range("G66").Value = 281100
$ws.Range("H66").Value = 228400
$ws.Range("I66").Value = 240500
$ws.Range("J66").Value = 223000
$ws.Range("D72").Value = 277800
$ws.Range("E72").Value = 308600
$ws.Range("F72").Value = 126500
$ws.Range("G72").Value = -3800
$ws.Range("H72").Value = -32200
$ws.Range("I72").Value = -94200
$ws.Range("J72").Value = -57200
$ws.Range("D76").Value = 525600
$ws.Range("E76").Value = 558300
$ws.Range("F76").Value = 353000
$ws.Range("G76").Value = 231200
$ws.Range("H76").Value = 204900
$ws.Range("I76").Value = 54400
$ws.Range("J76").Value = 38800
$ws.Range("D81").Value = -33700
$ws.Range("E81").Value = 182000
$ws.Range("F81").Value = 127400
$ws.Range("G81").Value = 31300
$ws.Range("H81").Value = 13900
$ws.Range("I81").Value = 11100
$ws.Range("J81").Value = -15200
$ws.Range("D83").Value = 28400
$ws.Range("E83").Value = 27800
$ws.Range("F83").Value = 23600
$ws.Range("G83").Value = 17700
$ws.Range("H83").Value = 15600
$ws.Range("I83").Value = 10900
$ws.Range("J83").Value = "NA"
$ws.Range("D89").Value = 51000
$ws.Range("E89").Value = 50600
$ws.Range("F89").Value = 158800
$ws.Range("G89").Value = 17300
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 25700
$ws.Range("J89").Value = -7600
$ws.Range("D91").Value = -130400
$ws.Range("E91").Value = -113700
$ws.Range("F91").Value = -72600
$ws.Range("G91").Value = -81400
$ws.Range("H91").Value = -57900
$ws.Range("I91").Value = -42900
$ws.Range("J91").Value = -61700
$ws.Range("D94").Value = -116900
$ws.Range("E94").Value = -1400
$ws.Range("F94").Value = -31100
$ws.Range("G94").Value = -61300
$ws.Range("H94").Value = -16100
$ws.Range("I94").Value = 9700
$ws.Range("J94").Value = "NA"
$ws.Range("D100").Value = 72300
$ws.Range("F100").Value = 3400
$ws.Range("G100").Value = 27100
$ws.Range("H100").Value = 96200
$ws.Range("I100").Value = -14100
$ws.Range("J100").Value = "NA"
$ws.Range("D101").Value = -12600
$ws.Range("E101").Value = 49300
$ws.Range("F101").Value = 18600
$ws.Range("G101").Value = 7300
$ws.Range("H101").Value = 2100
$ws.Range("I101").Value = -1600
$ws.Range("J101").Value = "NA"
$ws.Range("D102").Value = -6200
$ws.Range("E102").Value = 99000
$ws.Range("F102").Value = 149600
$ws.Range("G102").Value = -9600
$ws.Range("H102").Value = 84200
$ws.Range("I102").Value = 19700
$ws.Range("J102").Value = -12200
